$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.005.86"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "'2.237.61"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'304.90"
$ws.Range("E5").Value = "  -4.59%  "

$ws.Range("D6").Value = "'94.96"
$ws.Range("E6").Value = "  -7.16%  "

$ws.Range("D7").Value = "'0.568"
$ws.Range("E7").Value = "  -1.63%  "

$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  -6.15%  "

$ws.Range("D10").Value = "'34.51"
$ws.Range("E10").Value = "  -8.09%  "

$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = "  -3.41%  "

$ws.Range("D12").Value = "'7.17"
$ws.Range("E12").Value = "  -6.44%  "

$ws.Range("E13").Value = "  -2.86%  "

$ws.Range("D14").Value = "'2.576.66"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").Value = "'2.244.55"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").Value = "'0.815"
$ws.Range("E16").Value = "  -5.37%  "

$ws.Range("D17").Value = "'13.46"
$ws.Range("E17").Value = "  -7.27%  "

$ws.Range("D18").Value = "'43.833.28"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").Value = "'0.0₃0955"
$ws.Range("E19").Value = "  -3.42%  "

$ws.Range("D20").Value = "'12.12"
$ws.Range("E20").Value = "  -10.20%  "

$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "  -6.00%  "

$ws.Range("D22").Value = "'64.60"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").Value = "'237.55"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("E24").Value = "  -7.95%  "

$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("E26").Value = "  -8.60%  "

$ws.Range("D27").Value = "'9.84"
$ws.Range("E27").Value = "  -3.84%  "

$ws.Range("D28").Value = "'2.12"
$ws.Range("E28").Value = "  -3.55%  "

$ws.Range("D29").Value = "'36.33"
$ws.Range("E29").Value = "  -4.51%  "

$ws.Range("D30").Value = "'20.08"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("D31").Value = "'5.86"
$ws.Range("E31").Value = "  -6.14%  "

$ws.Range("D32").Value = "'152.82"
$ws.Range("E32").Value = "  -5.14%  "

$ws.Range("D33").Value = "'0.0804"
$ws.Range("E33").Value = "  -5.96%  "

$ws.Range("D34").Value = "'3.25"
$ws.Range("E34").Value = "  +6.48%  "

$ws.Range("D35").Value = "'2.63"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("D36").Value = "'0.110"
$ws.Range("E36").Value = "  -4.81%  "

$ws.Range("E37").Value = "  -1.11%  "

$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -9.66%  "

$ws.Range("D39").Value = "'14.82"
$ws.Range("E39").Value = "  -11.21%  "

$ws.Range("D40").Value = "'3.32"
$ws.Range("E40").Value = "  -11.54%  "

$ws.Range("D41").Value = "'3.78"
$ws.Range("E41").Value = "  -10.87%  "

$ws.Range("D42").Value = "'0.0299"
$ws.Range("E42").Value = "  -5.52%  "

$ws.Range("D44").Value = "'1.730.26"
$ws.Range("E44").Value = "  -3.92%  "

$ws.Range("D45").Value = "'84.87"
$ws.Range("E45").Value = "  +2.32%  "

$ws.Range("D46").Value = "'0.185"
$ws.Range("E46").Value = "  -7.31%  "

$ws.Range("D47").Value = "'99.66"
$ws.Range("E47").Value = "  -5.23%  "

$ws.Range("D48").Value = "'4.86"
$ws.Range("E48").Value = "  -7.13%  "

$ws.Range("D49").Value = "'14.53"
$ws.Range("E49").Value = "  +1.45%  "

$ws.Range("D50").Value = "'8.02"
$ws.Range("E50").Value = "  -4.12%  "

$ws.Range("D51").Value = "'68.33"
$ws.Range("E51").Value = "  -10.51%  "
